# Fix header spacing when unchecked
#
# The header currently reads (all on one visual line once the merge
# fields are filled in):
#   {{#SHOW_REPORT_UNDER}}{{REPORT_UNDER_TITLE}}{{REPORT_NUMBERS}}{{/SHOW_REPORT_UNDER}}
#
# When SHOW_REPORT_UNDER is shown, REPORT_UNDER_TITLE and REPORT_NUMBERS
# need to be on their own lines instead of running together, so we add a
# manual line break after the opening merge tag and another one after
# REPORT_UNDER_TITLE.

$d = $word.ActiveDocument

# The markers live in the primary (default) header of the (only) section.
$header = $d.Sections(1).Headers.Item(1)

# Make sure we are editing the paragraph we think we are before touching
# anything.
$searchRange = $header.Range.Duplicate
$found = $searchRange.Find.Execute("{{#SHOW_REPORT_UNDER}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the SHOW_REPORT_UNDER header marker"
}

$para = $header.Range.Paragraphs.Item(1)
$paraXml = $para.Range.WordOpenXML

# Pull the paragraph's own opening tag (keeps its paraId/rsid/etc. exactly
# as Word already has them) and its <w:pPr> (keeps the existing paragraph
# formatting - the right-aligned, bold/underlined header style - exactly
# as it is) instead of hard-coding them.
$pOpenStart = $paraXml.IndexOf("<w:p ")
$pOpenEnd = $paraXml.IndexOf(">", $pOpenStart) + 1
$pOpenTag = $paraXml.Substring($pOpenStart, $pOpenEnd - $pOpenStart)

$pPrStart = $paraXml.IndexOf("<w:pPr>")
$pPrEndTag = "</w:pPr>"
$pPrEnd = $paraXml.IndexOf($pPrEndTag) + $pPrEndTag.Length
$pPr = $paraXml.Substring($pPrStart, $pPrEnd - $pPrStart)

$runProps = "<w:rPr><w:b/><w:bCs/><w:u w:val=""single""/></w:rPr>"

$runs = ""
$runs += "<w:r>" + $runProps + "<w:t>{{#SHOW_REPORT_UNDER}}</w:t></w:r>"
$runs += "<w:r>" + $runProps + "<w:br/></w:r>"
$runs += "<w:r>" + $runProps + "<w:t>{{REPORT_UNDER_TITLE}}</w:t></w:r>"
$runs += "<w:r>" + $runProps + "<w:br/></w:r>"
$runs += "<w:r><w:t>{{REPORT_NUMBERS</w:t></w:r>"
$runs += "<w:proofErr w:type=""gramStart""/>"
$runs += "<w:r><w:t>}}{</w:t></w:r>"
$runs += "<w:proofErr w:type=""gramEnd""/>"
$runs += "<w:r><w:t>{/SHOW_REPORT_UNDER}}</w:t></w:r>"

$newParaInner = $pOpenTag + $pPr + $runs + "</w:p>"

$xmlPackage = "<pkg:package xmlns:pkg=""http://schemas.microsoft.com/office/2006/xmlPackage"">" +
    "<pkg:part pkg:name=""/word/document.xml"" pkg:contentType=""application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"">" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main"" xmlns:w14=""http://schemas.microsoft.com/office/word/2010/wordml"">" +
    "<w:body>" + $newParaInner + "</w:body></w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

# Replace the whole (single-paragraph) header range with the rebuilt
# paragraph - InsertXML replaces the contents of the range it is called on.
$header.Range.InsertXML($xmlPackage)
